function Get-ParagraphIndexByText($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $trimmed = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $text) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument

# Resolve the positions of the skill lines up front, while every line of text is
# still unique, so none of the lookups below are ambiguous.
$idxBases    = Get-ParagraphIndexByText $d "Bases de données : SQL, MongoDB, Neo4j, Redis"
$idxSoft     = Get-ParagraphIndexByText $d "Soft_Skills : flexible"
$idxVisu     = Get-ParagraphIndexByText $d "Visualisation : tableau"
$idxMLOps    = Get-ParagraphIndexByText $d "MLOps : software development, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"

# 1. "Bases de données : SQL, MongoDB, Neo4j, Redis" -> "Soft_Skills : flexible"
$d.Paragraphs.Item($idxBases).Range.Text = "Soft_Skills : flexible"

# 2. original "Soft_Skills : flexible" paragraph -> "Visualisation : tableau"
$d.Paragraphs.Item($idxSoft).Range.Text = "Visualisation : tableau"

# 3. insert a brand-new paragraph right after it for "MLOps : ..."
$d.Paragraphs.Item($idxSoft).Range.InsertParagraphAfter()
$d.Paragraphs.Item($idxSoft + 1).Range.Text = "MLOps : software development, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"

# The insertion shifted every paragraph after $idxSoft down by one slot, so the
# original "Visualisation : tableau" and "MLOps : ..." paragraphs are now one
# index further along than they used to be.
$idxVisu  = $idxVisu + 1
$idxMLOps = $idxMLOps + 1

# 4. remove the now-duplicated "Visualisation : tableau" paragraph
$d.Paragraphs.Item($idxVisu).Range.Delete()

# Deleting that paragraph shifts everything after it back up by one slot.
if ($idxMLOps -gt $idxVisu) {
    $idxMLOps = $idxMLOps - 1
}

# 5. original "MLOps : ..." -> "Bases de données : SQL, MongoDB, Neo4j, Redis"
$d.Paragraphs.Item($idxMLOps).Range.Text = "Bases de données : SQL, MongoDB, Neo4j, Redis"

Write-Output "done"
